# Update "Sprint 3" points-completed value for 2019-12-01 (row 40, column B)
# from 0 to 7. Column F (cumulative points completed) and column G (points
# left, = $E$3 - F) are driven by shared formulas, so Excel will recalc them
# automatically for rows 40-47 when B40 changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B40").Value = 7

# Fix the sidebar/view state: move the selection down to B41 and scroll the
# window so row 22 is the top visible row (matches the saved view position).
$ws.Range("B41").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
